$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen the sheet-tab area in the saved window view ---
try { $wb.Windows.Item(1).TabRatio = 0.5 } catch { }

# --- Insert a new first column (A), shifting the existing data (old A..H) to (B..I) ---
$ws.Columns("A").Insert()

# --- Give the new header cell (A1) the same look as the rest of the header row ---
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New "Subject Priority" column values ---
$ws.Range("A1").Value = "Subject Priority"
$ws.Range("A2").Value = "Algorithms/Data Structures"
$ws.Range("A3").Value = "Testing"
$ws.Range("A4").Value = "SQL/Xquery"
$ws.Range("A5").Value = "Version Control"
$ws.Range("A6").Value = "Design Patterns/Clean Code"
$ws.Range("A7").Value = "Functional Programming"
$ws.Range("A8").Value = "Scripting"
$ws.Range("A9").Value = "More Command Line Operations"

# --- Tidy up a stray trailing space in the pre-existing "linux" entry ---
$ws.Range("F2").Value = "linux"

# --- Give the new column its own width ---
$ws.Columns("A").ColumnWidth = 28.8697674418605

# --- Move the selection, matching the saved view state ---
$ws.Range("A17").Select() | Out-Null
